$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 30) mirroring the existing table's layout:
# data (date serial), hora (time string), preco (number), site (string)
$ws.Range("A30").Value = 45209
$ws.Range("A30").NumberFormat = $ws.Range("A29").NumberFormat

$ws.Range("B30").Value = "21:13"
$ws.Range("C30").Value = 95.90000000000001
$ws.Range("D30").Value = "natura"
